$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 225.09091
$ws.Range("I2").Value = 219.55556
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 219.55556
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -106.55556
$ws.Range("N2").Value = -476
$ws.Range("H61").Value = 731.63635
$ws.Range("I61").Value = 629.8
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 1889.4
$ws.Range("L61").Value = 5250
$ws.Range("M61").Value = -1717.4
$ws.Range("N61").Value = -5594
$ws.Range("H86").Value = 3105.8
$ws.Range("I86").Value = 3333.5
$ws.Range("K86").Value = 3333.5
$ws.Range("M86").Value = -2210.5
$ws.Range("H89").Value = 3105.8
$ws.Range("I89").Value = 3333.5
$ws.Range("K89").Value = 16667.5
$ws.Range("M89").Value = -11051.5
$ws.Range("H113").Value = 4939.2666
$ws.Range("I113").Value = 4339.2
$ws.Range("J113").Value = 6139.4
$ws.Range("K113").Value = 4339.2
$ws.Range("L113").Value = 6139.4
$ws.Range("M113").Value = -1085.2
$ws.Range("N113").Value = -12647.4
$ws.Range("H132").Value = 71433620
$ws.Range("I132").Value = 71433620
$ws.Range("K132").Value = 214300860
$ws.Range("M132").Value = -214298330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1155.5927
$ws.Range("I2").Value = 1005.7368
$ws.Range("J2").Value = 1511.5
$ws.Range("K2").Value = 1005.7368
$ws.Range("L2").Value = 1511.5
$ws.Range("M2").Value = -892.7368
$ws.Range("N2").Value = -1737.5
$ws.Range("H4").Value = 774
$ws.Range("I4").Value = 774
$ws.Range("K4").Value = 774
$ws.Range("M4").Value = -658
$ws.Range("H45").Value = 5080.5454
$ws.Range("I45").Value = 5088.6
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 5088.6
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -4711.6
$ws.Range("N45").Value = -5754
$ws.Range("H63").Value = 16000
$ws.Range("J63").Value = 16000
$ws.Range("L63").Value = 16000
$ws.Range("N63").Value = -17372
$ws.Range("H66").Value = 16000
$ws.Range("J66").Value = 16000
$ws.Range("L66").Value = 80000
$ws.Range("N66").Value = -86864
$ws.Range("H74").Value = 2533.4285
$ws.Range("I74").Value = 1998.6428
$ws.Range("K74").Value = 1998.6428
$ws.Range("M74").Value = -1124.6428
$ws.Range("H77").Value = 2533.4285
$ws.Range("I77").Value = 1998.6428
$ws.Range("K77").Value = 9993.214
$ws.Range("M77").Value = -5625.214
$ws.Range("H110").Value = 750
$ws.Range("I110").Value = 750
$ws.Range("K110").Value = 750
$ws.Range("M110").Value = 1295
$ws.Range("H116").Value = 1155.5927
$ws.Range("I116").Value = 1005.7368
$ws.Range("J116").Value = 1511.5
$ws.Range("K116").Value = 1005.7368
$ws.Range("L116").Value = 1511.5
$ws.Range("M116").Value = 1288.2632
$ws.Range("N116").Value = -6099.5
$ws.Range("H131").Value = 84900
$ws.Range("J131").Value = 84900
$ws.Range("L131").Value = 84900
$ws.Range("N131").Value = -94980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1155.5927
$ws.Range("I3").Value = 1005.7368
$ws.Range("J3").Value = 1511.5
$ws.Range("K3").Value = 1005.7368
$ws.Range("L3").Value = 1511.5
$ws.Range("M3").Value = -891.7368
$ws.Range("N3").Value = -1739.5
$ws.Range("H20").Value = 3271
$ws.Range("I20").Value = 2149.5
$ws.Range("K20").Value = 2149.5
$ws.Range("M20").Value = -1902.5
$ws.Range("H86").Value = 18687.125
$ws.Range("I86").Value = 21666.166
$ws.Range("K86").Value = 21666.166
$ws.Range("M86").Value = -20543.166
$ws.Range("H89").Value = 18687.125
$ws.Range("I89").Value = 21666.166
$ws.Range("K89").Value = 108330.83
$ws.Range("M89").Value = -102714.83
$ws.Range("H94").Value = 1044.8572
$ws.Range("I94").Value = 1044.8572
$ws.Range("K94").Value = 1044.8572
$ws.Range("M94").Value = -593.8571999999999
$ws.Range("H99").Value = 1234.2858
$ws.Range("I99").Value = 898.3333
$ws.Range("K99").Value = 898.3333
$ws.Range("M99").Value = 599.6667
$ws.Range("H105").Value = 2042.6
$ws.Range("I105").Value = 2348
$ws.Range("J105").Value = 1584.5
$ws.Range("K105").Value = 2348
$ws.Range("L105").Value = 1584.5
$ws.Range("M105").Value = -601
$ws.Range("N105").Value = -5078.5
$ws.Range("H130").Value = 34666.668
$ws.Range("J130").Value = 34666.668
$ws.Range("L130").Value = 34666.668
$ws.Range("N130").Value = -44706.668
$ws.Range("H134").Value = 1589.9
$ws.Range("I134").Value = 1712.625
$ws.Range("J134").Value = 1099
$ws.Range("K134").Value = 5137.875
$ws.Range("L134").Value = 3297
$ws.Range("M134").Value = -2602.875
$ws.Range("N134").Value = -8367

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 79937.5
$ws.Range("J20").Value = 79937.5
$ws.Range("L20").Value = 79937.5
$ws.Range("N20").Value = -80409.5
$ws.Range("H30").Value = 79937.5
$ws.Range("J30").Value = 79937.5
$ws.Range("L30").Value = 79937.5
$ws.Range("N30").Value = -80119.5
$ws.Range("H31").Value = 5745.227
$ws.Range("J31").Value = 3752.3333
$ws.Range("L31").Value = 3752.3333
$ws.Range("N31").Value = -4342.3333
$ws.Range("H34").Value = 5745.227
$ws.Range("J34").Value = 3752.3333
$ws.Range("L34").Value = 3752.3333
$ws.Range("N34").Value = -4156.3333
$ws.Range("H68").Value = 72857.07000000001
$ws.Range("J68").Value = 74615.30499999999
$ws.Range("L68").Value = 74615.30499999999
$ws.Range("N68").Value = -76113.30499999999
$ws.Range("H71").Value = 72857.07000000001
$ws.Range("J71").Value = 74615.30499999999
$ws.Range("L71").Value = 223845.915
$ws.Range("N71").Value = -231333.915
$ws.Range("H95").Value = 17204.6
$ws.Range("J95").Value = 17204.6
$ws.Range("L95").Value = 17204.6
$ws.Range("N95").Value = -22696.6
$ws.Range("H96").Value = 19883.111
$ws.Range("J96").Value = 19883.111
$ws.Range("L96").Value = 19883.111
$ws.Range("N96").Value = -25375.111
$ws.Range("H122").Value = 1572.1666
$ws.Range("I122").Value = 1211.3334
$ws.Range("K122").Value = 3634.0002
$ws.Range("M122").Value = -1184.0002
$ws.Range("H128").Value = 79937.5
$ws.Range("J128").Value = 79937.5
$ws.Range("L128").Value = 79937.5
$ws.Range("N128").Value = -89897.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 514.2857
$ws.Range("I17").Value = 312.5
$ws.Range("J17").Value = 783.3333
$ws.Range("K17").Value = 937.5
$ws.Range("L17").Value = 2349.9999
$ws.Range("M17").Value = -768.5
$ws.Range("N17").Value = -2687.9999
$ws.Range("H131").Value = 1998.9231
$ws.Range("I131").Value = 2495
$ws.Range("J131").Value = 1908.7273
$ws.Range("K131").Value = 7485
$ws.Range("L131").Value = 5726.1819
$ws.Range("M131").Value = -2445
$ws.Range("N131").Value = -15806.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 858.8182
$ws.Range("I97").Value = 483
$ws.Range("J97").Value = 2550
$ws.Range("K97").Value = 483
$ws.Range("L97").Value = 2550
$ws.Range("M97").Value = 13
$ws.Range("N97").Value = -3542
$ws.Range("H126").Value = 4498
$ws.Range("I126").Value = 4498
$ws.Range("K126").Value = 13494
$ws.Range("M126").Value = -11024
$ws.Range("H132").Value = 2529.5715
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1493.1666
$ws.Range("J22").Value = 1493.1666
$ws.Range("L22").Value = 1493.1666
$ws.Range("N22").Value = -2083.1666
$ws.Range("H27").Value = 1493.1666
$ws.Range("J27").Value = 1493.1666
$ws.Range("L27").Value = 1493.1666
$ws.Range("N27").Value = -1707.1666
$ws.Range("H55").Value = 1031
$ws.Range("I55").Value = 399.66666
$ws.Range("J55").Value = 1842.7142
$ws.Range("K55").Value = 399.66666
$ws.Range("L55").Value = 1842.7142
$ws.Range("M55").Value = -226.66666
$ws.Range("N55").Value = -2188.7142
$ws.Range("H100").Value = 1995.75
$ws.Range("I100").Value = 1995.75
$ws.Range("K100").Value = 1995.75
$ws.Range("M100").Value = -1454.75
$ws.Range("H128").Value = 69331.336
$ws.Range("J128").Value = 69331.336
$ws.Range("L128").Value = 69331.336
$ws.Range("N128").Value = -79291.336
$ws.Range("H136").Value = 4247
$ws.Range("J136").Value = 5994.5
$ws.Range("L136").Value = 17983.5
$ws.Range("N136").Value = -23083.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3994
$ws.Range("I62").Value = 3994
$ws.Range("K62").Value = 3994
$ws.Range("M62").Value = -3370
$ws.Range("H65").Value = 3994
$ws.Range("I65").Value = 3994
$ws.Range("K65").Value = 19970
$ws.Range("M65").Value = -16850
$ws.Range("H103").Value = 27499.5
$ws.Range("J103").Value = 27499.5
$ws.Range("L103").Value = 27499.5
$ws.Range("N103").Value = -29843.5
$ws.Range("H132").Value = 10578.75
$ws.Range("I132").Value = 10578.75
$ws.Range("K132").Value = 31736.25
$ws.Range("M132").Value = -29206.25
